# Atualização dos dados e melhorias no codigo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds data starting at row 1 (municipio rows, no
# header). We need to:
#   1. Insert two new rows at the top: a technical "Unnamed: N" header row
#      (bold, centered, top-aligned, boxed with a thin border) and a
#      human-readable header row ("municipio" / "Nº DE CASOS" /
#      "Óbitos confirmados").
#   2. Append two new data rows at the bottom of the table:
#      "outros paises" (27) and "outros estados" (14).

# --- 1. Shift existing data down by inserting two rows at the top -------
$ws.Rows.Item(1).Resize(2).Insert()

# --- 2. Row 1: technical header row --------------------------------------
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# --- 3. Row 2: human-readable header row ---------------------------------
$ws.Range("A2").Value = "municipio"
$ws.Range("B2").Value = "Nº DE CASOS"
$ws.Range("C2").Value = "Óbitos confirmados"

# --- 4. Append the two new trailing rows ---------------------------------
# (columns C41/C42 are left blank - no "Obitos confirmados" reported for
# these two aggregate rows, matching the rest of the table's empty cells)
$ws.Range("A41").Value = "outros paises"
$ws.Range("B41").Value = 27

$ws.Range("A42").Value = "outros estados"
$ws.Range("B42").Value = 14
